$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '67.296.57'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -1.24%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.310.17'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '186.80'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '578.06'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -1.03%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.40%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.129'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  -0.69%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '3.884.74'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("E13").Value = '  -0.48%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '27.45'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.38%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '67.545.01'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.89%  '
$ws.Range("E16").Value = '  -0.25%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.298.47'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.05%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '444.79'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +6.75%  '
$ws.Range("E19").Value = '  -0.80%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '13.59'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.85%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '7.73'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.47%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '73.92'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +3.36%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +1.85%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '3.454.95'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '
$ws.Range("E26").Value = '  +1.45%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.188'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.08'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -3.66%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("E30").Value = '  +1.38%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '22.90'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("E33").Value = '  -0.05%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.24'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("E36").Value = '  +4.68%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '162.96'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("E38").Value = '  -1.75%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '27.11'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("E41").Value = '  +0.15%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '2.759.17'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +3.75%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '6.25'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.96%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.0674'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '2.42'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '24.84'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.25%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '40.10'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -1.84%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '326.68'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.11%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.0274'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  +1.24%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '31.21'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +1.34%  '
